$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 190: 01-10-2021
$ws.Range("A190").Formula = '="01-10-2021"'
$ws.Range("A190").Copy()
$ws.Range("A190").PasteSpecial(-4163)
$ws.Range("B190").Value = -0.86
$ws.Range("C190").Value = -0.25
$ws.Range("D190").Value = 0.08

# Row 191: 04-10-2021
$ws.Range("A191").Formula = '="04-10-2021"'
$ws.Range("A191").Copy()
$ws.Range("A191").PasteSpecial(-4163)
$ws.Range("B191").Value = -0.86
$ws.Range("C191").Value = -0.24
$ws.Range("D191").Value = -0.05
